$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 3
$ws.Range("Q2").Value = 1.01

# Row 3
$ws.Range("F3").Value = 2.66
$ws.Range("J3").Value = 3.1
$ws.Range("K3").Value = 5.1
$ws.Range("P3").Value = 1.66
$ws.Range("Q3").Value = 1.9

# Row 4
$ws.Range("H4").Value = 3.15
$ws.Range("I4").Value = 4

# Row 6
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 4.6
$ws.Range("H6").Value = 1.81
$ws.Range("I6").Value = 1.95
$ws.Range("P6").Value = 2.38
$ws.Range("Q6").Value = 1.61

# Row 7
$ws.Range("F7").Value = 2.42
$ws.Range("H7").Value = 2.4
$ws.Range("J7").Value = 3.25
$ws.Range("K7").Value = 6.2
$ws.Range("P7").Value = 1.93
$ws.Range("Q7").Value = 1.66

# Row 8
$ws.Range("F8").Value = 3.75
$ws.Range("H8").Value = 1.9
$ws.Range("J8").Value = 3.45
$ws.Range("K8").Value = 3.85
$ws.Range("P8").Value = 1.73
$ws.Range("Q8").Value = 1.96

# Row 9
$ws.Range("F9").Value = 1.68
$ws.Range("H9").Value = 1.89

# Row 11
$ws.Range("F11").Value = 2.48
$ws.Range("G11").Value = 3.2
$ws.Range("H11").Value = 2.38
$ws.Range("I11").Value = 3.05
$ws.Range("J11").Value = 2.78

# Row 12
$ws.Range("G12").Value = 1.76
$ws.Range("P12").Value = 1.59

# Row 13
$ws.Range("G13").Value = 2.1
$ws.Range("J13").Value = 3.1
$ws.Range("Q13").Value = 1.6

# Row 14
$ws.Range("H14").Value = 3.95
$ws.Range("I14").Value = 4.5
$ws.Range("P14").Value = 2

# Row 15
$ws.Range("F15").Value = 3.8

# Row 16
$ws.Range("H16").Value = 1.86

# Row 17
$ws.Range("Q17").Value = 2.8
